$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B4").Value = 7.456999999999999
$ws.Range("A9").Value = -21.794
$ws.Range("B9").Value = 5.518
$ws.Range("C9").Value = -10.989
$ws.Range("B11").Value = 7.052
$ws.Range("A18").Value = -21.875
$ws.Range("A20").Value = -20.433
$ws.Range("B23").Value = 7.542
$ws.Range("B24").Value = 5.781000000000001
$ws.Range("B26").Value = 6.686999999999999
$ws.Range("A27").Value = -21.581
$ws.Range("C27").Value = -13.237
$ws.Range("C29").Value = -12.228
$ws.Range("C32").Value = -12.409
$ws.Range("B34").Value = 7.211999999999999
$ws.Range("A35").Value = -21.842
$ws.Range("B35").Value = 5.794000000000001
$ws.Range("C37").Value = -12.266
$ws.Range("C38").Value = -12.267
$ws.Range("C41").Value = -12.375
$ws.Range("C45").Value = -13.283
$ws.Range("B48").Value = 5.48
$ws.Range("B49").Value = 6.237
$ws.Range("C51").Value = -11.155
$ws.Range("B52").Value = 5.189000000000001
$ws.Range("C57").Value = -14.015
$ws.Range("C64").Value = -11.016
$ws.Range("B66").Value = 5.087000000000001
$ws.Range("B67").Value = 5.398999999999999
$ws.Range("A69").Value = -21.361
$ws.Range("A76").Value = -20.716
$ws.Range("A78").Value = -21.068
$ws.Range("B78").Value = 6.555000000000001
$ws.Range("B80").Value = 6.889999999999999
$ws.Range("A82").Value = -21.258
$ws.Range("C82").Value = -11.541
$ws.Range("A83").Value = -21.558
$ws.Range("A93").Value = -21.432
$ws.Range("C93").Value = -11.25
$ws.Range("B99").Value = 5.488
$ws.Range("C102").Value = -12.788
$ws.Range("B104").Value = 7.394000000000001
$ws.Range("C105").Value = -12.377
